# Converts the "Introdução à Eletrônica Digital" deck into the
# "Inteligência Artificial" deck:
#  - Slide 1: update title + subtitle text
#  - Slides 2-6: update title + bullet content
#  - Slides 7-10 (old): removed entirely

$p = $ppt.ActivePresentation

# --- Slide 1: title + subtitle -------------------------------------------
$s1 = $p.Slides.Item(1)

$t = $s1.Shapes.Item(1).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Inteligência Artificial"

$t = $s1.Shapes.Item(2).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Autor: Kézia"

# --- Slide 2: "O que é Inteligência Artificial?" --------------------------
$s2 = $p.Slides.Item(2)

$t = $s2.Shapes.Item(1).TextFrame.TextRange
$t.Text = "X"
$t.Text = "O que é Inteligência Artificial?"

$t = $s2.Shapes.Item(2).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Definição de Inteligência Artificial`rCapacidade de simular o pensamento humano através de algoritmos e sistemas computacionais"
[void]$t.Paragraphs(1,1).InsertBefore("`r")

# --- Slide 3: "Aplicações da Inteligência Artificial" ---------------------
$s3 = $p.Slides.Item(3)

$t = $s3.Shapes.Item(1).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Aplicações da Inteligência Artificial"

$t = $s3.Shapes.Item(2).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Setor de Saúde: diagnóstico médico, pesquisa de medicamentos`rSetor Automobilístico: carros autônomos`rSetor Financeiro: análise de risco, detecção de fraudes`rSetor de Varejo: recomendação de produtos, atendimento ao cliente"

# --- Slide 4: "Algoritmos de Inteligência Artificial" ---------------------
$s4 = $p.Slides.Item(4)

$t = $s4.Shapes.Item(1).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Algoritmos de Inteligência Artificial"

$t = $s4.Shapes.Item(2).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Aprendizado de Máquina (Machine Learning)`rRedes Neurais Artificiais`rProcessamento de Linguagem Natural`rVisão Computacional"

# --- Slide 5: "Benefícios e Desafios" --------------------------------------
$s5 = $p.Slides.Item(5)

$t = $s5.Shapes.Item(1).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Benefícios e Desafios"

$t = $s5.Shapes.Item(2).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Benefícios: automação de tarefas, aumento da eficiência, criação de soluções inovadoras`rDesafios: viés nas decisões, falta de transparência, privacidade e segurança dos dados"

# --- Slide 6: "Impacto da Inteligência Artificial no Futuro" --------------
$s6 = $p.Slides.Item(6)

$t = $s6.Shapes.Item(1).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Impacto da Inteligência Artificial no Futuro"

$t = $s6.Shapes.Item(2).TextFrame.TextRange
$t.Text = "X"
$t.Text = "Transformação de diversos setores da economia`rNovos empregos e necessidade de reskilling`rAdaptação da sociedade à presença da IA"

# --- Remove the old slides 7-10 (Circuitos Combinacionais, Circuitos
#     Sequenciais, Microcontroladores, Conclusão) -------------------------
$p.Slides.Item(10).Delete()
$p.Slides.Item(9).Delete()
$p.Slides.Item(8).Delete()
$p.Slides.Item(7).Delete()

Write-Output "done: $($p.Slides.Count) slides remain"
